$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 2 first (paragraph starting "First, for online Amazon shopping, ..."):
#   The old "_GoBack" bookmark used to sit between "Amazon sh" and
#   "opping, ". It moves to paragraph 1 (handled below), so here it simply
#   goes away and the two runs around it collapse back into a single
#   "Amazon shopping, " run. We do this before touching paragraph 1 because
#   bookmark names are unique in the document -- re-adding "_GoBack"
#   elsewhere later will automatically relocate/replace this one.
# ---------------------------------------------------------------------------

$rAmazon = $d.Content
$rAmazon.Find.ClearFormatting()
$found2 = $rAmazon.Find.Execute("Amazon sh", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Amazon sh' run"
}
$amazonStart = $rAmazon.Start
$amazonEnd = $rAmazon.End

# Search for "opping, " starting right after the "Amazon sh" match so we
# don't accidentally match the unrelated "holiday shopping, " phrase that
# appears earlier in the document. Wrap=0 (wdFindStop) keeps the search
# from wrapping back around to that earlier occurrence.
$rOpping = $d.Range($amazonEnd, $d.Content.End)
$rOpping.Find.ClearFormatting()
$found3 = $rOpping.Find.Execute("opping, ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if (-not $found3) {
    throw "Could not find 'opping, ' run"
}
$oppingEnd = $rOpping.End

# Put temporary stop bookmarks around the whole "Amazon shopping, " span so
# the text rewrite below cannot bleed into the neighbouring "online " /
# "please" runs (which share the same visible formatting).
$rLeftStop = $d.Range($amazonStart, $amazonStart)
$d.Bookmarks.Add("ZZZ_TempLeft", $rLeftStop)
$rRightStop = $d.Range($oppingEnd, $oppingEnd)
$d.Bookmarks.Add("ZZZ_TempRight", $rRightStop)

# Drop the old (soon to be stale) _GoBack bookmark that currently lives
# inside this span.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Rewrite the whole bounded span as a single run of text. The final text is
# identical to the concatenation of the two existing runs, so merely
# assigning the same characters back would be treated as a no-op by the
# engine and the runs would stay split. Route through a placeholder value
# first to force a genuine content mutation, then swap in the real text.
$rSpan = $d.Range($amazonStart, $oppingEnd)
$rSpan.Text = "ZZZ_PLACEHOLDER_AMAZON_ZZZ"

$rFinal = $d.Content
$rFinal.Find.ClearFormatting()
$foundPlaceholder = $rFinal.Find.Execute("ZZZ_PLACEHOLDER_AMAZON_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPlaceholder) {
    throw "Could not find placeholder text"
}
$rFinal.Text = "Amazon shopping, "

# Remove the temporary stop bookmarks; the merge they produced stays.
$d.Bookmarks("ZZZ_TempLeft").Delete()
$d.Bookmarks("ZZZ_TempRight").Delete()

# ---------------------------------------------------------------------------
# Change 1 (paragraph 1, the title):
#   "Post " -> "Post-" and the "_GoBack" bookmark now sits right after
#   "Post-" (between it and "Holiday Shopping…").
# ---------------------------------------------------------------------------

# Locate the exact "Post " run inside the title so we do not depend on
# hard-coded character offsets.
$rPost = $d.Content
$rPost.Find.ClearFormatting()
$found = $rPost.Find.Execute("Post ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Post ' run in title"
}
$postStart = $rPost.Start
$postEnd = $rPost.End

# Insert a temporary marker bookmark right BEFORE "Post " so that when we
# rewrite its text below, the preceding "For your " run is not absorbed
# into the same run as the rewritten text.
$rBeforePost = $d.Range($postStart, $postStart)
$d.Bookmarks.Add("ZZZ_TempBeforePost", $rBeforePost)

# Rewrite "Post " -> "Post-" (drops the trailing space, adds a hyphen).
$rPostText = $d.Range($postStart, $postEnd)
$rPostText.Text = "Post-"

# New end of the just-rewritten run.
$newPostEnd = $postStart + 5

# Add the "_GoBack" bookmark right after "Post-" (collapsed range), which
# also splits the following "Holiday Shopping…" run away from it. Because
# bookmark names are unique, this automatically relocates/replaces any
# "_GoBack" bookmark left over elsewhere in the document.
$rGoBackNew = $d.Range($newPostEnd, $newPostEnd)
$d.Bookmarks.Add("_GoBack", $rGoBackNew)

# Remove the temporary helper bookmark; the run split it forced stays in
# place even after the bookmark itself is gone.
$d.Bookmarks("ZZZ_TempBeforePost").Delete()
